# DOMA-1858 — add new column "Дата передачи показаний" (Reading submission
# date) between "Показание 4" (J) and "Дата поверки" (old K, now L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at K; everything from K.. shifts right to L..
$ws.Columns("K").Insert()

# Give the new column the same formatting as its neighbour (the old "Дата
# поверки" column, now shifted to L) - border/fill/font match the rest of
# the date columns.
$ws.Range("L1:L10").Copy()
$ws.Range("K1:K10").PasteSpecial(-4122)

# Header
$ws.Range("K1").Value = "Дата передачи показаний"

# Data - one submission date per address block
$ws.Range("K2:K6").Value = "2021-12-20"
$ws.Range("K7:K10").Value = "2021-12-21"

# Column width matching the source workbook
$ws.Columns("K").ColumnWidth = 20
